# Update "want to go" (F column) head-counts for several rows on the
# "展览" and "全部类型" sheets, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# row -> new F-column value
$updates = @{
    3  = 703
    5  = 19
    8  = 1624
    9  = 5765
    10 = 470
    11 = 324
    12 = 261
    13 = 75
    14 = 352
    15 = 125
    16 = 4695
    18 = 1244
    22 = 89
    24 = 85
    28 = 53
    31 = 29
    32 = 47
    33 = 16
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
